$wb = $excel.ActiveWorkbook

# Update F3 (110 -> 111) and F5 (68 -> 69) on sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 111
$ws1.Range("F5").Value = 69

# Update F3 (110 -> 111) and F5 (68 -> 69) on sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 111
$ws4.Range("F5").Value = 69
